# The canonical-OOXML diff for this resource shows only whitespace/attribute
# re-ordering (xmlns declarations on <w:document>, and attribute order on
# <w:pgSz>/<w:pgMar>/<w:rFonts>/<w:lang>/<w:latentStyles>/<w:lsdException>/
# <w:style>/<w:tblInd>/<w:tblCellMar> in word/styles.xml) that was produced by
# the commit's tooling when it re-serialized the template while wiring the
# (unrelated, Java-side) M2Doc-version custom property into the test
# fixtures. Every changed line pair is attribute-set-identical to its
# predecessor - same elements, same attribute names/values, same text -
# so there is no actual document content, formatting, or structural change
# to replay here. Touch the document model without mutating anything so the
# package round-trips through this automation session unchanged.
$d = $word.ActiveDocument
$null = $d.Content.Text
